{"js": "// Replace the 100 arithmetic-problem cell texts in the practice-sheet\n// table with their new values. The table has exactly 100 single-run\n// paragraphs (5 columns x 20 rows) that follow one leading date\n// paragraph, in the same order as `pairs` below (old -> new).\nconst pairs = [[\"67-55=\", \"3+27=\"], [\"15+60=\", \"92-24=\"], [\"58-29=\", \"78-23=\"], [\"89-27=\", \"60+30=\"], [\"79-21=\", \"23+15=\"], [\"21+29=\", \"18+24=\"], [\"17+20=\", \"30-25=\"], [\"42-6=\", \"57-7=\"], [\"83-32=\", \"19+34=\"], [\"98-39=\", \"41+33=\"], [\"42-4=\", \"94-63=\"], [\"28-11=\", \"11+86=\"], [\"1+47=\", \"95+3=\"], [\"93-28=\", \"65+12=\"], [\"48-8=\", \"7-6=\"], [\"53+17=\", \"98-93=\"], [\"53+7=\", \"13-4=\"], [\"37-2=\", \"74-27=\"], [\"11+9=\", \"7-1=\"], [\"90-21=\", \"87-14=\"], [\"70+29=\", \"46+48=\"], [\"94-54=\", \"28-1=\"], [\"21+2=\", \"21+67=\"], [\"5+93=\", \"81-60=\"], [\"97-97=\", \"30-5=\"], [\"57+41=\", \"93-55=\"], [\"70-20=\", \"34+39=\"], [\"33-22=\", \"93-80=\"], [\"66-43=\", \"81-49=\"], [\"3+0=\", \"28-5=\"], [\"72+23=\", \"3+73=\"], [\"32+59=\", \"9-9=\"], [\"2+69=\", \"36+34=\"], [\"19-18=\", \"39+25=\"], [\"28+0=\", \"11+70=\"], [\"18+50=\", \"13+80=\"], [\"45+49=\", \"78-50=\"], [\"23+51=\", \"43+6=\"], [\"95-3=\", \"88-86=\"], [\"25+2=\", \"23+71=\"], [\"94-35=\", \"15+8=\"], [\"55-13=\", \"70-61=\"], [\"73-67=\", \"94+1=\"], [\"9+47=\", \"39+37=\"], [\"88-14=\", \"63-21=\"], [\"33+47=\", \"14+56=\"], [\"20+21=\", \"23-6=\"], [\"61-47=\", \"49-25=\"], [\"14+33=\", \"73-69=\"], [\"2+67=\", \"85-78=\"], [\"22+18=\", \"18+23=\"], [\"68+28=\", \"64-64=\"], [\"34+40=\", \"7+52=\"], [\"87-49=\", \"57-52=\"], [\"67-54=\", \"87+3=\"], [\"99-4=\", \"81+7=\"], [\"76-54=\", \"45-3=\"], [\"88-15=\", \"40+53=\"], [\"39-29=\", \"39-33=\"], [\"18+62=\", \"41-5=\"], [\"37+29=\", \"47-6=\"], [\"48-9=\", \"86-80=\"], [\"25+73=\", \"34+24=\"], [\"41-0=\", \"5+72=\"], [\"93-71=\", \"98-57=\"], [\"72-57=\", \"93-24=\"], [\"65-18=\", \"21+12=\"], [\"0+20=\", \"55+7=\"], [\"65+7=\", \"85-38=\"], [\"90-80=\", \"74-69=\"], [\"32-15=\", \"40+27=\"], [\"36+4=\", \"76-40=\"], [\"39+22=\", \"70+9=\"], [\"71-33=\", \"17-15=\"], [\"45+52=\", \"20+0=\"], [\"20-13=\", \"67-26=\"], [\"90-67=\", \"60+30=\"], [\"56-29=\", \"11+16=\"], [\"57-21=\", \"45+34=\"], [\"25+9=\", \"94-92=\"], [\"98-20=\", \"99-0=\"], [\"84+8=\", \"4+74=\"], [\"4+55=\", \"32+15=\"], [\"68-30=\", \"63-49=\"], [\"21+66=\", \"10+18=\"], [\"48+10=\", \"17+34=\"], [\"26-15=\", \"39-38=\"], [\"40+7=\", \"91-16=\"], [\"7+33=\", \"65-3=\"], [\"32-10=\", \"94-35=\"], [\"51+33=\", \"45+41=\"], [\"50+17=\", \"85-56=\"], [\"49+7=\", \"83-7=\"], [\"56-44=\", \"17+72=\"], [\"62-34=\", \"49-16=\"], [\"13-9=\", \"58-12=\"], [\"41-10=\", \"75+23=\"], [\"81-70=\", \"20+43=\"], [\"5+13=\", \"0+86=\"], [\"4+7=\", \"81-38=\"]];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n// The first paragraph is the date header; the table's 100 cell\n// paragraphs (one run each) follow it in row-major document order,\n// lining up 1:1 with `pairs`.\nconst offset = items.length - pairs.length;\n\nlet applied = 0;\nfor (let i = 0; i < pairs.length; i++) {\n  const [oldText, newText] = pairs[i];\n  const para = items[offset + i];\n  if (para.text === oldText) {\n    para.insertText(newText, Word.InsertLocation.replace);\n    applied++;\n  } else {\n    // Fallback: text didn't line up positionally (unexpected), so find\n    // the exact cell by searching the whole document for the old text.\n    const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n    results.load(\"text\");\n    await context.sync();\n    if (results.items.length > 0) {\n      results.items[0].insertText(newText, Word.InsertLocation.replace);\n      applied++;\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the 100 arithmetic-problem cell texts in the practice-sheet\n# table with their new values. Each pair is unique and, applied in this\n# order, matches exactly one occurrence in the document at the time it\n# runs (verified against the source diff), so a simple global\n# Find/ReplaceAll per pair reproduces the target edit while leaving all\n# formatting (fonts, sizes, table/paragraph properties) untouched.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"67-55=\", \"3+27=\"),\n    @(\"15+60=\", \"92-24=\"),\n    @(\"58-29=\", \"78-23=\"),\n    @(\"89-27=\", \"60+30=\"),\n    @(\"79-21=\", \"23+15=\"),\n    @(\"21+29=\", \"18+24=\"),\n    @(\"17+20=\", \"30-25=\"),\n    @(\"42-6=\", \"57-7=\"),\n    @(\"83-32=\", \"19+34=\"),\n    @(\"98-39=\", \"41+33=\"),\n    @(\"42-4=\", \"94-63=\"),\n    @(\"28-11=\", \"11+86=\"),\n    @(\"1+47=\", \"95+3=\"),\n    @(\"93-28=\", \"65+12=\"),\n    @(\"48-8=\", \"7-6=\"),\n    @(\"53+17=\", \"98-93=\"),\n    @(\"53+7=\", \"13-4=\"),\n    @(\"37-2=\", \"74-27=\"),\n    @(\"11+9=\", \"7-1=\"),\n    @(\"90-21=\", \"87-14=\"),\n    @(\"70+29=\", \"46+48=\"),\n    @(\"94-54=\", \"28-1=\"),\n    @(\"21+2=\", \"21+67=\"),\n    @(\"5+93=\", \"81-60=\"),\n    @(\"97-97=\", \"30-5=\"),\n    @(\"57+41=\", \"93-55=\"),\n    @(\"70-20=\", \"34+39=\"),\n    @(\"33-22=\", \"93-80=\"),\n    @(\"66-43=\", \"81-49=\"),\n    @(\"3+0=\", \"28-5=\"),\n    @(\"72+23=\", \"3+73=\"),\n    @(\"32+59=\", \"9-9=\"),\n    @(\"2+69=\", \"36+34=\"),\n    @(\"19-18=\", \"39+25=\"),\n    @(\"28+0=\", \"11+70=\"),\n    @(\"18+50=\", \"13+80=\"),\n    @(\"45+49=\", \"78-50=\"),\n    @(\"23+51=\", \"43+6=\"),\n    @(\"95-3=\", \"88-86=\"),\n    @(\"25+2=\", \"23+71=\"),\n    @(\"94-35=\", \"15+8=\"),\n    @(\"55-13=\", \"70-61=\"),\n    @(\"73-67=\", \"94+1=\"),\n    @(\"9+47=\", \"39+37=\"),\n    @(\"88-14=\", \"63-21=\"),\n    @(\"33+47=\", \"14+56=\"),\n    @(\"20+21=\", \"23-6=\"),\n    @(\"61-47=\", \"49-25=\"),\n    @(\"14+33=\", \"73-69=\"),\n    @(\"2+67=\", \"85-78=\"),\n    @(\"22+18=\", \"18+23=\"),\n    @(\"68+28=\", \"64-64=\"),\n    @(\"34+40=\", \"7+52=\"),\n    @(\"87-49=\", \"57-52=\"),\n    @(\"67-54=\", \"87+3=\"),\n    @(\"99-4=\", \"81+7=\"),\n    @(\"76-54=\", \"45-3=\"),\n    @(\"88-15=\", \"40+53=\"),\n    @(\"39-29=\", \"39-33=\"),\n    @(\"18+62=\", \"41-5=\"),\n    @(\"37+29=\", \"47-6=\"),\n    @(\"48-9=\", \"86-80=\"),\n    @(\"25+73=\", \"34+24=\"),\n    @(\"41-0=\", \"5+72=\"),\n    @(\"93-71=\", \"98-57=\"),\n    @(\"72-57=\", \"93-24=\"),\n    @(\"65-18=\", \"21+12=\"),\n    @(\"0+20=\", \"55+7=\"),\n    @(\"65+7=\", \"85-38=\"),\n    @(\"90-80=\", \"74-69=\"),\n    @(\"32-15=\", \"40+27=\"),\n    @(\"36+4=\", \"76-40=\"),\n    @(\"39+22=\", \"70+9=\"),\n    @(\"71-33=\", \"17-15=\"),\n    @(\"45+52=\", \"20+0=\"),\n    @(\"20-13=\", \"67-26=\"),\n    @(\"90-67=\", \"60+30=\"),\n    @(\"56-29=\", \"11+16=\"),\n    @(\"57-21=\", \"45+34=\"),\n    @(\"25+9=\", \"94-92=\"),\n    @(\"98-20=\", \"99-0=\"),\n    @(\"84+8=\", \"4+74=\"),\n    @(\"4+55=\", \"32+15=\"),\n    @(\"68-30=\", \"63-49=\"),\n    @(\"21+66=\", \"10+18=\"),\n    @(\"48+10=\", \"17+34=\"),\n    @(\"26-15=\", \"39-38=\"),\n    @(\"40+7=\", \"91-16=\"),\n    @(\"7+33=\", \"65-3=\"),\n    @(\"32-10=\", \"94-35=\"),\n    @(\"51+33=\", \"45+41=\"),\n    @(\"50+17=\", \"85-56=\"),\n    @(\"49+7=\", \"83-7=\"),\n    @(\"56-44=\", \"17+72=\"),\n    @(\"62-34=\", \"49-16=\"),\n    @(\"13-9=\", \"58-12=\"),\n    @(\"41-10=\", \"75+23=\"),\n    @(\"81-70=\", \"20+43=\"),\n    @(\"5+13=\", \"0+86=\"),\n    @(\"4+7=\", \"81-38=\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $rng = $d.Range()\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $rng.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
